# ProjectAllocations.xlsx -> rebuilt as a "University applications" sheet.
# The previous "preallocation" report (7 cols x 4 rows, hyperlinked names,
# bold+underline styling) is replaced wholesale with a simple 5-col x 5-row
# University/Projects/Applied/Accepted table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three name hyperlinks from the old report before wiping content.
$ws.Hyperlinks.Delete()

# Wipe all existing cell content/formatting.
$ws.Cells.Clear()

# Drop the old 7-column layout so we can rebuild clean A:E column metadata.
$ws.Columns("A:G").Delete()

# Restore default gridline display (round-tripping can otherwise flip this).
$excel.ActiveWindow.DisplayGridlines = $true

# --- Header row ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "University"
$ws.Range("C1").Value = "Projects"
$ws.Range("D1").Value = "Applied"
$ws.Range("E1").Value = "Accepted"

# --- Data rows ---
$ws.Range("A2").Value = "Jose Brown"
$ws.Range("B2").Value = "UoG"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

$ws.Range("A3").Value = "Darion Considine"
$ws.Range("B3").Value = "UoG"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("A4").Value = "Penelope Nienow"
$ws.Range("B4").Value = "UoG"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = "Judy Parker"
$ws.Range("B5").Value = "UoG"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Bold header row, matching the workbook's existing 12pt header style.
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A1:E1").Font.Size = 12

# Size the five columns to fit their new content.
$ws.Columns("A:E").AutoFit()

# Land the selection on E1, like the source workbook.
$ws.Range("E1").Select() | Out-Null
